# Add files via upload
# Populates columns O:T (Minute3a, Second3b, Rep3Row, rank, points, Team)
# on the "Score" sheet for rows 2-21, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

# O, P, Q, T raw values per row (Minute3a, Second3b, Rep3Row, Team)
$data = @{
    2  = @(11, 15, 75, 1602)
    3  = @(10, 59, 66, 1602)
    4  = @(10, 37, 92, 1602)
    5  = @(10, 59, 102, 1602)
    6  = @(10, 39, 76, 1602)
    7  = @(10, 57, 97, 1602)
    8  = @(12, 10, 80, 1602)
    9  = @(13, 7, 108, 1602)
    10 = @(10, 5, 89, 1602)
    11 = @(10, 56, 73, 1602)
    12 = @(13, 8, 134, 1602)
    13 = @(10, 40, 75, 1602)
    14 = @(13, 10, 100, 1602)
    15 = @(11, 27, 97, 1602)
    16 = @(11, 57, 100, 1602)
    17 = @(11, 26, 91, 1602)
    18 = @(10, 30, 68, 1602)
    19 = @(11, 40, 101, 1602)
    20 = @(10, 12, 97, 1602)
    21 = @(14, 54, 61, 1602)
}

foreach ($row in 2..21) {
    $vals = $data[$row]

    $ws.Cells.Item($row, 15).Value = $vals[0]   # O - Minute3a
    $ws.Cells.Item($row, 16).Value = $vals[1]   # P - Second3b
    $ws.Cells.Item($row, 17).Value = $vals[2]   # Q - Rep3Row

    $ws.Cells.Item($row, 18).Formula = "=INT((O$row*60+P$row-Q$row)/60)"   # R - rank
    $ws.Cells.Item($row, 19).Formula = "=O$row*60+P$row-Q$row-R$row*60"    # S - points

    $ws.Cells.Item($row, 20).Value = $vals[3]   # T - Team
}

# Update the active selection shown in the saved file
$ws.Range("O22").Select()
